$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B21").Value = 18
$ws.Range("C21").Value = "Create a dBase of only tags"
$ws.Range("C21").WrapText = $true

$ws.Range("C22").Select()
